# Updated cryptos list on Mon Aug 26 03:09:45 UTC 2024 with GitHub Actions
# Applies the latest price / volume(1h) snapshot values to the cryptocurrency
# table on the active worksheet, including the NEARProtocol / EthereumClassic
# row swap (rows 33-34).
#
# Numeric-looking "Price" values are entered with a leading apostrophe so
# Excel keeps them as text (matching the source data, which stores prices
# like "569.48" or "20.00" as literal strings rather than numbers); the
# Style reset immediately afterwards strips the implicit "quote prefix"
# formatting that the apostrophe would otherwise leave behind, so cell
# formatting stays exactly as it was before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.025.59"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.736.64"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'569.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").Value = "'158.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.46%  "
$ws.Range("D11").Value = "'5.71"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.56%  "
$ws.Range("D12").Value = "'0.382"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "3.218.76"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'26.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "63.633.70"
$ws.Range("E15").Value = "  -0.52%  "
$ws.Range("E16").Value = "  -2.05%  "
$ws.Range("D17").Value = "2.741.33"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'12.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "'4.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'354.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.64%  "
$ws.Range("D21").Value = "'6.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "'0.522"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.21%  "
$ws.Range("D24").Value = "'64.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("D25").Value = "'0.171"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'8.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("D29").Value = "'1.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").Value = "'7.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  +7.18%  "
$ws.Range("D32").Value = "'163.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'4.90"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'20.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D37").Value = "'1.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "'0.989"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "'350.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.87%  "
$ws.Range("D40").Value = "'6.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("E41").Value = "  -1.39%  "
$ws.Range("D42").Value = "'38.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").Value = "'22.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "'21.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.21%  "
$ws.Range("D45").Value = "'0.0582"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").Value = "'134.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").Value = "'0.624"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
